$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Mệnh Không Thân Kiếp"
$ws.Range("B3").Value = "Bạn là người khôn ngoan, sắc sảo nên ông trời thử thách bạn với những hoàn cảnh trớ trêu."
$ws.Range("C3").Value = "Bạn cũng là tuýp người nhiệt tình 5 phút, cả thèm mau chán."

$ws.Range("J8").Select()
